$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the data of data rows 2 and 3 (species records), while
# leaving the columns that already held identical values in both rows
# untouched (C, I, K, N, P..AE, AG..AY).

# ---- Row 2: becomes the former row-3 record (Trådticka) ----
$ws.Range("A2").Value = 111813166
$ws.Range("B2").Value = 90087
$ws.Range("D2").Value = "LC"
$ws.Range("E2").Value = 3298
$ws.Range("F2").Value = "Trådticka"
$ws.Range("G2").Value = "Climacocystis borealis"
$ws.Range("H2").Value = "(Fr.) Kotl. & Pouzar"

# J2 did not exist before; it now exists as an (empty) text cell.
$ws.Range("J2").Value = "'"
$ws.Range("J2").ClearFormats()

# L2 and M2 existed before (L2 empty, M2 = "färska spår"); they no
# longer exist afterwards.
$ws.Range("L2").ClearContents()
$ws.Range("M2").ClearContents()

# AF2 did not exist before; it now exists as an (empty) text cell.
$ws.Range("AF2").Value = "'"
$ws.Range("AF2").ClearFormats()

# ---- Row 3: becomes the former row-2 record (Tretåig hackspett) ----
$ws.Range("A3").Value = 111813153
$ws.Range("B3").Value = 56398
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 100109
$ws.Range("F3").Value = "Tretåig hackspett"
$ws.Range("G3").Value = "Picoides tridactylus"
$ws.Range("H3").Value = "(Linnaeus, 1758)"

# J3 existed before (empty); it no longer exists afterwards.
$ws.Range("J3").ClearContents()

# L3 did not exist before; it now exists as an (empty) text cell.
$ws.Range("L3").Value = "'"
$ws.Range("L3").ClearFormats()

# M3 did not exist before; it now exists holding "färska spår".
$ws.Range("M3").Value = "färska spår"

# AF3 existed before (empty); it no longer exists afterwards.
$ws.Range("AF3").ClearContents()
